$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $rng = $Worksheet.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '57.202.17'
Set-TextValue $ws 'E2' '  +5.27%  '
Set-TextValue $ws 'D3' '3.245.97'
Set-TextValue $ws 'E3' '  +2.36%  '
Set-TextValue $ws 'E4' '  +0.13%  '
Set-TextValue $ws 'D5' '394.94'
Set-TextValue $ws 'E5' '  -0.85%  '
Set-TextValue $ws 'D6' '107.64'
Set-TextValue $ws 'E6' '  -1.39%  '
Set-TextValue $ws 'D7' '0.592'
Set-TextValue $ws 'E7' '  +8.21%  '
Set-TextValue $ws 'D8' '3.243.71'
Set-TextValue $ws 'E8' '  +2.50%  '
Set-TextValue $ws 'E9' '  +0.07%  '
Set-TextValue $ws 'D10' '0.624'
Set-TextValue $ws 'E10' '  +1.29%  '
Set-TextValue $ws 'D11' '39.05'
Set-TextValue $ws 'E11' '  +0.33%  '
Set-TextValue $ws 'D12' '0.0983'
Set-TextValue $ws 'E12' '  +11.87%  '
Set-TextValue $ws 'E13' '  +1.73%  '
Set-TextValue $ws 'D14' '3.759.04'
Set-TextValue $ws 'E14' '  +2.60%  '
Set-TextValue $ws 'D15' '8.17'
Set-TextValue $ws 'E15' '  +1.71%  '
Set-TextValue $ws 'D16' '19.06'
Set-TextValue $ws 'E16' '  -0.44%  '
Set-TextValue $ws 'D17' '3.240.27'
Set-TextValue $ws 'E17' '  +2.47%  '
Set-TextValue $ws 'D18' '1.02'
Set-TextValue $ws 'E18' '  -3.38%  '
Set-TextValue $ws 'D19' '10.75'
Set-TextValue $ws 'E19' '  +2.69%  '
Set-TextValue $ws 'D20' '56.964.66'
Set-TextValue $ws 'E20' '  +5.04%  '
Set-TextValue $ws 'D21' '3.34'
Set-TextValue $ws 'E21' '  +1.90%  '
Set-TextValue $ws 'D22' '0.0000112'
Set-TextValue $ws 'E22' '  +14.27%  '
Set-TextValue $ws 'D23' '12.89'
Set-TextValue $ws 'E23' '  +0.51%  '
Set-TextValue $ws 'D24' '295.14'
Set-TextValue $ws 'E24' '  +8.54%  '
Set-TextValue $ws 'E25' '  +4.45%  '
Set-TextValue $ws 'D26' '3.14'
Set-TextValue $ws 'E26' '  -3.48%  '
Set-TextValue $ws 'D27' '27.82'
Set-TextValue $ws 'E27' '  +0.67%  '
Set-TextValue $ws 'D28' '7.57'
Set-TextValue $ws 'E28' '  -5.09%  '
Set-TextValue $ws 'D29' '7.21'
Set-TextValue $ws 'E29' '  -2.18%  '
Set-TextValue $ws 'D30' '0.168'
Set-TextValue $ws 'E30' '  -1.07%  '
Set-TextValue $ws 'E31' '  -0.05%  '
Set-TextValue $ws 'D32' '11.32'
Set-TextValue $ws 'E32' '  +3.09%  '
Set-TextValue $ws 'D33' '0.107'
Set-TextValue $ws 'E33' '  -4.02%  '
Set-TextValue $ws 'D34' '39.01'
Set-TextValue $ws 'E34' '  +5.81%  '
Set-TextValue $ws 'D35' '0.0479'
Set-TextValue $ws 'E35' '  -4.43%  '
Set-TextValue $ws 'E36' '  +1.94%  '
Set-TextValue $ws 'D37' '51.62'
Set-TextValue $ws 'E37' '  +2.20%  '
Set-TextValue $ws 'D38' '0.998'
Set-TextValue $ws 'E38' '  -0.09%  '
Set-TextValue $ws 'D39' '3.47'
Set-TextValue $ws 'E39' '  -4.54%  '
Set-TextValue $ws 'E40' '  +2.10%  '
Set-TextValue $ws 'D41' '134.82'
Set-TextValue $ws 'E41' '  +3.53%  '
Set-TextValue $ws 'E42' '  +4.14%  '
Set-TextValue $ws 'D43' '17.04'
Set-TextValue $ws 'E43' '  -1.62%  '
Set-TextValue $ws 'E44' '  -1.22%  '
Set-TextValue $ws 'D45' '3.92'
Set-TextValue $ws 'E45' '  -4.66%  '
Set-TextValue $ws 'D46' '0.280'
Set-TextValue $ws 'E46' '  -3.93%  '
Set-TextValue $ws 'D47' '22.12'
Set-TextValue $ws 'E47' '  -0.49%  '
Set-TextValue $ws 'E48' '  +3.14%  '
Set-TextValue $ws 'D49' '2.156.29'
Set-TextValue $ws 'E49' '  +3.38%  '
Set-TextValue $ws 'B50' 'ApeXProtocol'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws 'D50' '2.34'
Set-TextValue $ws 'E50' '  -3.87%  '
Set-TextValue $ws 'B51' 'ThetaToken'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws 'D51' '1.95'
Set-TextValue $ws 'E51' '  +16.70%  '
